$d = $word.ActiveDocument

# Replace the two "date de création" / "date de version" fields: 13/07/2018 -> 16/07/2018
$d.Content.Find.Execute("13/07/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "16/07/2018", 2)
$d.Content.Find.Execute("13/07/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "16/07/2018", 2)

# Feature table price estimate: 1230€ -> 2785€
$d.Content.Find.Execute("1230€", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2785€", 2)

# Total estimate (before taxes figure): 24470 -> 26025
$d.Content.Find.Execute("24470", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "26025", 2)

# Final total amount: 50330 -> 51885
$d.Content.Find.Execute("50330", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "51885", 2)
